$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.284.07"
$ws.Range("E2").Value = "  -3.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.969.46"
$ws.Range("E3").Value = "  -4.14%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.15"
$ws.Range("E5").Value = "  -3.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.79"
$ws.Range("E7").Value = "  -10.07%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.06"
$ws.Range("E10").Value = "  -4.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +4.07%  "

$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.857"
$ws.Range("E13").Value = "  -7.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.34"
$ws.Range("E14").Value = "  +8.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.01"
$ws.Range("E15").Value = "  -8.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.256.80"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.42"
$ws.Range("E17").Value = "  -2.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.969.70"
$ws.Range("E18").Value = "  -4.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.183.04"
$ws.Range("E19").Value = "  -3.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.17"
$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  -2.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.45"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  -5.68%  "

$ws.Range("E26").Value = "  -5.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.79"
$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.88"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.84"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("E30").Value = "  +6.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.87"
$ws.Range("E32").Value = "  -7.03%  "

$ws.Range("E33").Value = "  -6.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0618"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -7.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("E36").Value = "  +2.52%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.28"
$ws.Range("E38").Value = "  -6.28%  "

$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +7.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0991"
$ws.Range("E41").Value = "  -2.71%  "

$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("E44").Value = "  -3.13%  "

$ws.Range("E45").Value = "  -4.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.30"
$ws.Range("E46").Value = "  -3.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.94"
$ws.Range("E47").Value = "  -6.69%  "

$ws.Range("E48").Value = "  -6.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.341.59"
$ws.Range("E49").Value = "  -5.55%  "

$ws.Range("E50").Value = "  -3.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.152.24"
